# Junction_Flooding_477 — replace the 4 data rows with a new simulation run
# (custom accuracy + 1000 data points), drop the now-unused 6th row, and
# nudge a handful of column widths back to their "default" 7/8-char size.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New data for rows 2-5 (Time + J1..J33), one row per flooding timestep.
# ---------------------------------------------------------------------------
$row2 = @(45158.50694444445,22.58,15.542,4.221,47.493,39.284,17.769,58.8,27.341,11.61,17.881,18.828,19.728,5.673,17.67,24.849,14.79,3.779,2.46,261.617,49.202,16.31,32.642,17.025,2.109,29.266,14.407,12.944,15.145,19.485,3.64,51.902,9.071,20.391)
$row3 = @(45158.51388888889,23.06,16.637,2.04,49.52,40.967,18.147,70.22499999999999,27.923,12.226,18.412,19.904,20.858,5.797,18.046,25.579,15.268,1.68,1.257,267.375,50.484,16.657,33.744,17.925,2.208,34.006,14.713,13.168,15.447,20.795,1.294,63.516,9.329000000000001,20.825)
$row4 = @(45158.52083333334,1.922,0.991,0.838,3.761,3.126,1.517,14.52,2.327,0.985,1.309,1.592,1.497,0.515,1.504,2.189,1.524,0.961,0.383,15.713,4.784,1.388,3.094,1.687,0.07000000000000001,6.144,1.226,1.279,1.456,1.579,0.784,13.614,0.655,1.746)
$row5 = @(45158.52777777778,24.02,17.71,1.36,51.94,42.96,18.9,68.54000000000001,29.09,12.9,19.31,20.94,21.98,6.04,18.8,26.73,15.81,0.9,0.93,278.8,52.44,17.35,35.27,18.74,2.36,33.89,15.33,13.59,15.97,21.94,0.5600000000000001,61.76,9.779999999999999,21.69)

$newRows = @($row2, $row3, $row4, $row5)

for ($r = 0; $r -lt $newRows.Length; $r++) {
    $vals = $newRows[$r]
    $excelRow = $r + 2
    $arr = New-Object 'object[,]' 1,$vals.Length
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $arr[0,$i] = $vals[$i]
    }
    $rng = $ws.Range("A" + $excelRow + ":AH" + $excelRow)
    $rng.Value2 = $arr
}

# ---------------------------------------------------------------------------
# 2. Row 6 no longer exists in the new dataset (only 4 data rows remain).
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).Delete()

# ---------------------------------------------------------------------------
# 3. A few columns flip between the "7" and "8" character custom widths.
#    Excel's COM ColumnWidth is expressed in the default font's character
#    units and is offset from the raw OOXML <col width> by ~0.8333 (5px of
#    padding at 96dpi / Calibri 11), so width=7 -> ColumnWidth=6.1666...7
#    and width=8 -> ColumnWidth=7.1666...7.
# ---------------------------------------------------------------------------
$widthSeven  = 6.1666666666666667
$widthEight  = 7.1666666666666667

$colsToSeven  = @(2)
$colsToEight  = @(3,7,10,11,12,13,15,17,22,24,27,28,29,30,34)

foreach ($c in $colsToSeven) {
    $ws.Columns.Item($c).ColumnWidth = $widthSeven
}
foreach ($c in $colsToEight) {
    $ws.Columns.Item($c).ColumnWidth = $widthEight
}
